# Applies the "Append: 2025-09-08 12:36 JST" update to the
# "ランサーズ" (Lancers) sheet:
#   - Every existing row's timestamp (column A) is refreshed to the new
#     scrape time.
#   - Two brand-new listings are inserted into the ranked list (at what
#     become rows 3 and 7), pushing the remaining rows down.
#   - Hyperlinks in column F are rebuilt so each one keeps pointing at
#     the correct row after the inserts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-09-08 12:36:36"

# --- 1. Insert the two new rows, pushing everything below them down ---
# Insert at row 3 first (new listing about the in-house logistics system).
$ws.Rows.Item(3).Insert()
# Insert at row 7 (new listing about the admin-system detailed design).
# This happens *after* the first insert, so row numbers below already
# reflect the first shift.
$ws.Rows.Item(7).Insert()

# --- 2. Write out the full row data (rows 2-13) in final form ---

# Row 2 (unchanged listing, only timestamp refreshed)
$ws.Range("A2").Value = $newTimestamp

# Row 3 (NEW)
$ws.Range("A3").Value = $newTimestamp
$ws.Range("B3").Value = "自社開発のロジシステムをサポート及び開発できる方募集【PHP, Python, VBA etc】"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5389460"
$ws.Range("G3").Value = 305
$ws.Range("H3").Value = "🔥Python ◆開発 ○PHP"

# Row 4 (was row 3)
$ws.Range("A4").Value = $newTimestamp

# Row 5 (was row 4)
$ws.Range("A5").Value = $newTimestamp

# Row 6 (was row 5)
$ws.Range("A6").Value = $newTimestamp

# Row 7 (NEW)
$ws.Range("A7").Value = $newTimestamp
$ws.Range("B7").Value = "【急募】管理システムの詳細設計・開発を依頼します!"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5389414"
$ws.Range("G7").Value = 115
$ws.Range("H7").Value = "◆開発 ◇管理"

# Row 8 (was row 6)
$ws.Range("A8").Value = $newTimestamp

# Row 9 (was row 7)
$ws.Range("A9").Value = $newTimestamp

# Row 10 (was row 8)
$ws.Range("A10").Value = $newTimestamp

# Row 11 (was row 9)
$ws.Range("A11").Value = $newTimestamp

# Row 12 (was row 10)
$ws.Range("A12").Value = $newTimestamp

# Row 13 (was row 11)
$ws.Range("A13").Value = $newTimestamp

# --- 3. Rebuild hyperlinks in column F so each ref lines up correctly ---
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5389316")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5389460")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5273634")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5314730")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5388877")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5389414")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5389306")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5388922")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5372984")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5385681")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5389241")
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5389081")

# --- 4. Make sure the sheet's used-range dimension matches (A1:H13) ---
$ws.Range("A1").Select()
